$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("storageAssets")

# Row 14: lossFactor_WpK value (column G) changes from 1 to 0.8
$ws.Range("G14").Value = 0.8

# Row 16: newly re-added 10MWh Grid battery storage asset
$ws.Range("A16").Value = 13
$ws.Range("B16").Value = "Grid_battery_10MWh"
$ws.Range("C16").Value = "STORAGE"
$ws.Range("D16").Value = "STORAGE_ELECTRIC"
$ws.Range("E16").Value = 2000
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0.8
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = 0
$ws.Range("M16").NumberFormat = "0.00E+00"
$ws.Range("O16").Value = 0

# Update the view to reflect where the user ended up after the edit
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("N16").Select()
